# Update the 25 "two-digit x two-digit" problems in the single table.
# Each populated row (1, 5, 10, 15, 20) holds 5 expressions, one per column.
# Cells are addressed positionally (Table.Cell(row, col)) rather than via a
# global text Find/Replace, because a couple of the new values coincide with
# other cells' old values (e.g. row15/col3 becomes "64x89=", which is also
# the *old* value of row20/col2) - a sequential text replace could clobber
# itself. Direct cell addressing sets every value exactly once, unambiguously.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "33×25="; New = "81×37=" },
    @{ Row = 1;  Col = 2; Old = "73×46="; New = "85×45=" },
    @{ Row = 1;  Col = 3; Old = "47×57="; New = "26×61=" },
    @{ Row = 1;  Col = 4; Old = "84×60="; New = "35×91=" },
    @{ Row = 1;  Col = 5; Old = "87×56="; New = "39×88=" },

    @{ Row = 5;  Col = 1; Old = "85×50="; New = "51×55=" },
    @{ Row = 5;  Col = 2; Old = "61×36="; New = "76×94=" },
    @{ Row = 5;  Col = 3; Old = "37×95="; New = "39×16=" },
    @{ Row = 5;  Col = 4; Old = "62×97="; New = "35×36=" },
    @{ Row = 5;  Col = 5; Old = "45×17="; New = "91×34=" },

    @{ Row = 10; Col = 1; Old = "27×77="; New = "48×26=" },
    @{ Row = 10; Col = 2; Old = "20×19="; New = "18×70=" },
    @{ Row = 10; Col = 3; Old = "82×74="; New = "52×78=" },
    @{ Row = 10; Col = 4; Old = "46×80="; New = "92×30=" },
    @{ Row = 10; Col = 5; Old = "79×32="; New = "51×39=" },

    @{ Row = 15; Col = 1; Old = "18×20="; New = "43×67=" },
    @{ Row = 15; Col = 2; Old = "49×80="; New = "95×74=" },
    @{ Row = 15; Col = 3; Old = "33×34="; New = "64×89=" },
    @{ Row = 15; Col = 4; Old = "21×41="; New = "26×90=" },
    @{ Row = 15; Col = 5; Old = "36×41="; New = "96×99=" },

    @{ Row = 20; Col = 1; Old = "57×64="; New = "96×23=" },
    @{ Row = 20; Col = 2; Old = "64×89="; New = "21×79=" },
    @{ Row = 20; Col = 3; Old = "59×52="; New = "72×73=" },
    @{ Row = 20; Col = 4; Old = "20×93="; New = "69×94=" },
    @{ Row = 20; Col = 5; Old = "73×68="; New = "89×52=" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $u.Old) {
        Write-Host "WARNING: R$($u.Row)C$($u.Col) expected '$($u.Old)' but found '$current'"
    }
    $cell.Range.Text = $u.New
}

Write-Host "Done: updated $($updates.Count) expressions."
